# SanProspero report update: insert the missing 2021-02-08 data point, which
# shifts every subsequent 7-day-rolling-sum (col C) and per-100k rate (col D)
# value down by one row, and append two new trailing days
# (2021-03-01 / 2021-03-02).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Row 91 (2021-02-06) rolling sum now only covers 7 days that no
#        longer straddle the later gap the same way -> recomputed values.
$ws.Cells(91,3).Value = 3
$ws.Cells(91,4).Value = 49.90850108135086

# --- 2. Insert the missing day (serial 44235 = 2021-02-08) as new row 93,
#        pushing the former rows 93..113 down to 94..114.
$ws.Rows.Item(93).Insert()

$ws.Cells(93,1).Value = 44235
$ws.Cells(93,2).Value = 0
$ws.Cells(93,3).Value = 3
$ws.Cells(93,4).Value = 49.90850108135086

# Match the date-column formatting used by every other row in col A.
$ws.Cells(93,1).NumberFormat = $ws.Cells(92,1).NumberFormat()
$ws.Cells(93,1).Borders.LineStyle = 1
$ws.Cells(93,1).HorizontalAlignment = -4108
$ws.Cells(93,1).VerticalAlignment = -4160

# --- 3. Recomputed rolling-sum (C) / rate (D) values for the rows that got
#        shifted down by the insert (now at rows 94..112).
$updates = @(
    @{R=94;  C=4;  D=66.54466810846782},
    @{R=95;  C=6;  D=99.81700216270171},
    @{R=96;  C=6;  D=99.81700216270171},
    @{R=97;  C=7;  D=116.4531691898187},
    @{R=98;  C=11; D=182.9978372982865},
    @{R=99;  C=11; D=182.9978372982865},
    @{R=100; C=11; D=182.9978372982865},
    @{R=101; C=11; D=182.9978372982865},
    @{R=102; C=8;  D=133.0893362169356},
    @{R=103; C=9;  D=149.7255032440526},
    @{R=104; C=13; D=216.2701713525204},
    @{R=105; C=10; D=166.3616702711695},
    @{R=106; C=9;  D=149.7255032440526},
    @{R=107; C=13; D=216.2701713525204},
    @{R=108; C=16; D=266.1786724338713},
    @{R=109; C=21; D=349.359507569456},
    @{R=110; C=21; D=349.359507569456},
    @{R=111; C=22; D=365.9956745965729},
    @{R=112; C=22; D=365.9956745965729}
)

foreach ($u in $updates) {
    $ws.Cells($u.R, 3).Value = $u.C
    $ws.Cells($u.R, 4).Value = $u.D
}

# --- 4. Append a brand-new trailing row for 2021-03-02 (serial 44257),
#        with 1 new positive case and no rolling-sum data yet.
$newRow = 115

# Clone formatting from the row directly above (same date-column look).
$ws.Cells($newRow-1,1).Copy()
$ws.Cells($newRow,1).PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells($newRow,1).Value = 44257
$ws.Cells($newRow,2).Value = 1

Write-Output "edit applied"
